# perbaikan Gate dan radio button form
# Replace the leftover template placeholders for the signer's role and
# name with plain instructional placeholder text, in both occurrences
# of the signature block.

$d = $word.ActiveDocument

# "${roleInspektur" + "}" (split across two runs) -> single run with
# instructional text, e.g. "(Masukkan jabatan penandatangan)".
$d.Content.Find.Execute('${roleInspektur}', $true, $false, $false, $false, $false, $true, 1, $false, '(Masukkan jabatan penandatangan)', 2)

# "${inspektur}" -> "(Masukkan nama penandatangan)"
$d.Content.Find.Execute('${inspektur}', $true, $false, $false, $false, $false, $true, 1, $false, '(Masukkan nama penandatangan)', 2)
